$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Use an out-of-the-way helper cell to enter "0.015" as text (forcing Text
# number format there), then move it into E2 and restore the header-row
# formatting (style already used by B1/C1/D1) via PasteSpecial so we don't
# leave a stray custom number-format style behind.
$helper = $ws.Range("Z100")
$helper.NumberFormat = "@"
$helper.Value = "0.015"
$helper.Copy()
$ws.Range("E2").PasteSpecial(-4163)
$ws.Range("B1").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$helper.Clear()

# Fill the rest of the new "P value" column (E) with "<0.001"
$ws.Range("E3").Value = "<0.001"
$ws.Range("E4").Value = "<0.001"
$ws.Range("E5").Value = "<0.001"
$ws.Range("E6").Value = "<0.001"
$ws.Range("E7").Value = "<0.001"
$ws.Range("E8").Value = "<0.001"
$ws.Range("E9").Value = "<0.001"
$ws.Range("E10").Value = "<0.001"
$ws.Range("E11").Value = "<0.001"

# E3 keeps the centered body-row style already used in column A-D (style used by B3)
$ws.Range("E3").HorizontalAlignment = -4108

# Update the selection / scroll position to match the final workbook view
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("E11").Select()
